# Insert a new data row at row 359, shifting all existing rows (359-447)
# down by one (to 360-448), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("359:359").Insert()

$ws.Range("A359").Value = 6
$ws.Range("B359").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C359").Value = "Metropolitana"
$ws.Range("D359").Value = 44722
$ws.Range("E359").Value = 13
$ws.Range("F359").Value = 100112039
$ws.Range("G359").Value = "Ciboulette"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 760
$ws.Range("K359").Value = 700
$ws.Range("L359").Value = 800
$ws.Range("M359").Value = 758
$ws.Range("N359").Value = "`$/docena de atados"
$ws.Range("O359").Value = "Región Metropolitana"
$ws.Range("P359").Value = 253
$ws.Range("Q359").Value = 3
$ws.Range("R359").Value = "Hortaliza"
